$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 26740
$ws.Range("F4").Value = 590
$ws.Range("F6").Value = 614
$ws.Range("F7").Value = 176
$ws.Range("F8").Value = 553
$ws.Range("F9").Value = 232
$ws.Range("F15").Value = 71
$ws.Range("F16").Value = 431
$ws.Range("F17").Value = 60
$ws.Range("F18").Value = 1550
$ws.Range("F20").Value = 49

$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 201
$ws.Range("F7").Value = 201
$ws.Range("F10").Value = 441
$ws.Range("F15").Value = 62
$ws.Range("F16").Value = 25

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 5083
$ws.Range("F3").Value = 237

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 5083
$ws.Range("F4").Value = 237
$ws.Range("F5").Value = 26740
$ws.Range("F6").Value = 590
$ws.Range("F10").Value = 614
$ws.Range("F13").Value = 176
$ws.Range("F14").Value = 201
$ws.Range("F15").Value = 201
$ws.Range("F18").Value = 441
$ws.Range("F19").Value = 553
$ws.Range("F21").Value = 232
$ws.Range("F28").Value = 71
$ws.Range("F31").Value = 431
$ws.Range("F32").Value = 60
$ws.Range("F33").Value = 62
$ws.Range("F34").Value = 1551
$ws.Range("F36").Value = 25
$ws.Range("F37").Value = 49
